$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings are not
# auto-converted to Excel numbers (the source data is plain text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "46.238.22"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").Value = "2.608.35"
$ws.Range("E3").Value = "  +3.68%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "307.54"
$ws.Range("E5").Value = "  +2.74%  "

$ws.Range("D6").Value = "100.20"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "0.603"
$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  +5.60%  "

$ws.Range("D10").Value = "39.51"
$ws.Range("E10").Value = "  +5.70%  "

$ws.Range("E11").Value = "  +4.37%  "

$ws.Range("D12").Value = "54.15"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("E13").Value = "  +5.45%  "

$ws.Range("D14").Value = "3.006.75"
$ws.Range("E14").Value = "  +3.54%  "

$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").Value = "2.600.26"
$ws.Range("E16").Value = "  +3.51%  "

$ws.Range("D17").Value = "0.919"
$ws.Range("E17").Value = "  +4.38%  "

$ws.Range("E18").Value = "  +1.96%  "

$ws.Range("D19").Value = "46.425.57"
$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("E20").Value = "  +3.95%  "

$ws.Range("D21").Value = "12.92"
$ws.Range("E21").Value = "  -4.37%  "

$ws.Range("D22").Value = "6.73"
$ws.Range("E22").Value = "  +3.30%  "

$ws.Range("D23").Value = "71.40"
$ws.Range("E23").Value = "  +3.53%  "

$ws.Range("D24").Value = "272.94"
$ws.Range("E24").Value = "  +9.01%  "

$ws.Range("E25").Value = "  +5.07%  "

$ws.Range("E26").Value = "  +4.76%  "

$ws.Range("D27").Value = "29.18"
$ws.Range("E27").Value = "  +26.53%  "

$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "4.02"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").Value = "10.58"
$ws.Range("E30").Value = "  +4.34%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.27"
$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "38.88"
$ws.Range("E32").Value = "  -4.91%  "

$ws.Range("D33").Value = "6.34"
$ws.Range("E33").Value = "  +8.42%  "

$ws.Range("E34").Value = "  -4.37%  "

$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("E36").Value = "  +3.13%  "

$ws.Range("E37").Value = "  +3.53%  "

$ws.Range("D38").Value = "150.78"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").Value = "0.121"
$ws.Range("E39").Value = "  +4.30%  "

$ws.Range("E40").Value = "  +4.02%  "

$ws.Range("D41").Value = "23.19"
$ws.Range("E41").Value = "  +36.84%  "

$ws.Range("D42").Value = "15.90"
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("E43").Value = "  +6.69%  "

$ws.Range("D44").Value = "0.0331"
$ws.Range("E44").Value = "  +4.91%  "

$ws.Range("D45").Value = "4.05"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").Value = "2.116.61"
$ws.Range("E46").Value = "  +5.71%  "

$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").Value = "93.22"
$ws.Range("E48").Value = "  +1.98%  "

$ws.Range("E49").Value = "  +7.98%  "

$ws.Range("D50").Value = "1.78"
$ws.Range("E50").Value = "  -1.38%  "

$ws.Range("D51").Value = "108.90"
$ws.Range("E51").Value = "  +2.17%  "

# Restore default (Normal) style on column D so the text format applied
# above does not leave a lingering explicit cell style.
$ws.Range("D2:D51").Style = "Normal"